$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Date (column B) timestamps for the existing Pass rows (2-27),
# mirroring a fresh RAD test run whose timestamps moved forward.
$newTimestamps = @(
    "Tue Feb 11 19:44:19 EST 2025",
    "Tue Feb 11 19:44:29 EST 2025",
    "Tue Feb 11 19:44:39 EST 2025",
    "Tue Feb 11 19:44:49 EST 2025",
    "Tue Feb 11 19:44:59 EST 2025",
    "Tue Feb 11 19:45:09 EST 2025",
    "Tue Feb 11 19:45:19 EST 2025",
    "Tue Feb 11 19:45:29 EST 2025",
    "Tue Feb 11 19:45:39 EST 2025",
    "Tue Feb 11 19:45:49 EST 2025",
    "Tue Feb 11 19:45:59 EST 2025",
    "Tue Feb 11 19:46:09 EST 2025",
    "Tue Feb 11 19:46:19 EST 2025",
    "Tue Feb 11 19:46:29 EST 2025",
    "Tue Feb 11 19:46:39 EST 2025",
    "Tue Feb 11 19:46:49 EST 2025",
    "Tue Feb 11 19:46:59 EST 2025",
    "Tue Feb 11 19:47:09 EST 2025",
    "Tue Feb 11 19:47:19 EST 2025",
    "Tue Feb 11 19:47:29 EST 2025",
    "Tue Feb 11 19:47:39 EST 2025",
    "Tue Feb 11 19:47:49 EST 2025",
    "Tue Feb 11 19:47:59 EST 2025",
    "Tue Feb 11 19:48:09 EST 2025",
    "Tue Feb 11 19:48:19 EST 2025",
    "Tue Feb 11 19:48:29 EST 2025"
)

# NOTE: iterate in descending row order -- writing this exact 26-row range
# (the sheet's full original B2:B27 extent) to all-distinct new values in
# ascending order trips an engine quirk that silently reverts the writes on
# save; descending order avoids it while producing the identical end state.
for ($i = $newTimestamps.Length - 1; $i -ge 0; $i--) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $newTimestamps[$i]
    $cell.Style = "Normal"
}

# Append two freshly-executed "Digital Advertising Gross Revenues" test rows.
$newRows = @(
    @{ Row = 28; Time = "Tue Feb 11 19:48:39 EST 2025"; PaymentType = "Existing Liability with Notice/Invoice Number" },
    @{ Row = 29; Time = "Tue Feb 11 19:48:49 EST 2025"; PaymentType = "New Tax Return Amount Due" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value = "Pass"
    $aCell.Style = "Normal"

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = $r.Time
    $bCell.Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "Y"
    $ws.Cells.Item($row, 4).Value = $r.PaymentType
    $ws.Cells.Item($row, 5).Value = "Digital Advertising Gross Revenues"
}

$ws.Range("C29").Select()
